$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.611.66'
$ws.Range("E2").Value = '  -1.24%  '

$ws.Range("D3").Value = '1.844.00'
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4239'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3641'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07280'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8894'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.79%  '

$ws.Range("D13").Value = '1.832.31'
$ws.Range("E13").Value = '  -2.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.568'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.345'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06877'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008878'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9994'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("E21").Value = '  -3.13%  '

$ws.Range("D22").Value = '27.603.85'
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.987'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.45%  '

$ws.Range("E24").Value = '  -4.39%  '

$ws.Range("D25").Value = '2.031.95'
$ws.Range("E25").Value = '  -3.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.934'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '122.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.285'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.880'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08925'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7698'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.576'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.923'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.096'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9992'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  -1.22%  '

$ws.Range("E39").Value = '  -1.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01932'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.814'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.79%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.879'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.84%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5085'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.27%  '

$ws.Range("E44").Value = '  -2.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.274'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06580'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.16%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4734'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.15%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9992'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("E51").Value = '  -2.59%  '
